# Highlight quantitative metrics (percentages, dollar amounts, large numbers)
# within specific bullet paragraphs by splitting the plain run into
# bold + colored (2C3E50) runs around each metric token, leaving the
# surrounding text as separate plain runs - matching the hybrid
# bold+color highlighting described in the commit message.

$d = $word.ActiveDocument

# RGB(0x2C,0x3E,0x50) expressed as the BGR-packed long that Word's
# Font.Color / wdColor property expects (R + G*256 + B*65536).
$metricColor = 5258796

function Apply-MetricHighlights {
    param($ParaExact, $Metrics)

    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $fullText = $p.Range.Text.TrimEnd()
        if ($fullText -ne $ParaExact) {
            continue
        }

        $scan = $p.Range
        $paraEnd = $scan.End
        $cursor = $scan.Start

        foreach ($metric in $Metrics) {
            $search = $d.Range($cursor, $paraEnd)
            $found = $search.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
            if ($found) {
                $search.Font.Bold = 1
                $search.Font.Color = $metricColor
                $cursor = $search.End
            }
        }

        break
    }
}

$bullet = [char]0x2022
$plusMinus = [char]0x00B1

Apply-MetricHighlights ($bullet + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%") @("23%", "64%")

Apply-MetricHighlights ($bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $plusMinus + "4.2% to " + $plusMinus + "2.1%") @("87%", "71%", ($plusMinus + "4.2%"), ($plusMinus + "2.1%"))

Apply-MetricHighlights ($bullet + " Wrote RFP and analyzed bids from 1,200 vendors for research platform development") @("1,200")

Apply-MetricHighlights ($bullet + ' Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+') @('$400M', '$1B')

Apply-MetricHighlights ($bullet + ' Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M') @("73.5%", '$4.7M')

Apply-MetricHighlights ($bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%") @("87%", "71%")
